# "try something new in fsl_cleaning and add iycf_flags"
#
# This script adds a new "IYCF" worksheet (flag_name / Rationale / Action
# table for IYCF data-quality flags) at the end of the workbook, mirroring
# the structure of the existing FSL/WASH/NUT/MORT flag-description sheets.

$wb = $excel.ActiveWorkbook

# Add the new sheet right after the current last sheet (MORT)
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "IYCF"

# Header row
$ws.Range("A1").Value = " flag_name"
$ws.Range("B1").Value = "Rationale"
$ws.Range("C1").Value = "Action"

# flag_yes_foods
$ws.Range("A2").Value = "flag_yes_foods"
$ws.Range("B2").Value = "Respondent reported child consuming all foods groups"
$ws.Range("C2").Value = "Review with enumerators since its unlikely and implausible that the child consumed ALL the food groups"

# flag_yes_liquids
$ws.Range("A3").Value = "flag_yes_liquids"
$ws.Range("B3").Value = "Respondent reported child consuming all liquids groups"
$ws.Range("C3").Value = "Review with enumerators since its unlikely and implausible that the child consumed ALL the liquid groups"

# flag_no_anything
$ws.Range("A4").Value = "flag_no_anything"
$ws.Range("B4").Value = "Respondent reported child consuming no foods or liquids groups at all"
$ws.Range("C4").Value = "Check child age if <9 months and IF exclusively breastfeeding. If child is not exclusively breastfeeding and >9 months, then review with enumerator as its extremely unlikely that child did not consume any foods or liquids."

# flag_no_foods
$ws.Range("A5").Value = "flag_no_foods"
$ws.Range("B5").Value = "Respondent reported child consuming no foods groups while reporting eating solid or semi-solid food meals"
$ws.Range("C5").Value = "Review with enumerators since its erroneous that the child was reported to consume NO food groups YET ate solid or semi-solid food meals"

# flag_all_foods_no_meal
$ws.Range("A6").Value = "flag_all_foods_no_meal"
$ws.Range("B6").Value = "Respondent reported child consuming all foods groups while reporting not eating any solid or semi-solid food meals"
$ws.Range("C6").Value = "Review with enumerators since its erroneous that the child was reported to consume ALL 8 food groups YET did not eat any solid or semi-solid food meals"

# flag_some_foods_no_meal
$ws.Range("A7").Value = "flag_some_foods_no_meal"
$ws.Range("B7").Value = "Respondent reported child consuming some foods groups while reporting not eating any solid or semi-solid food meals"
$ws.Range("C7").Value = "Review with enumerators since its erroneous that the child was reported to consume SOME food groups YET did not eat any solid or semi-solid food meals"

# flag_high_mdd_low_mmf
$ws.Range("A8").Value = "flag_high_mdd_low_mmf"
$ws.Range("B8").Value = "Respondent reported high mdd score while reporting low meal frequency consumed (<=1)"
$ws.Range("C8").Value = "If MDD is high (>4) and meal frequency is low (<=1), review for enumerator bias/error"

# flag_under6_nobf_nomilk
$ws.Range("A9").Value = "flag_under6_nobf_nomilk"
$ws.Range("B9").Value = "Respondent reported child under 6 month and not breastfed and no milk given"
$ws.Range("C9").Value = "If child is <6months AND not breastfed AND no other milks given, review for enumerator bias/error"

# flag_meats_nostaples
$ws.Range("A10").Value = "flag_meats_nostaples"
$ws.Range("B10").Value = "Respondent reported child consuming meats but no staples"
$ws.Range("C10").Value = "Review with enumerators since its unlikely and implausible that the child consumed meats but NO staples."

# Column widths: A ~= 25.29 (bestFit-like), B = 15
$ws.Columns.Item(1).ColumnWidth = 24.43
$ws.Columns.Item(2).ColumnWidth = 14.2

# Make IYCF the active sheet/tab with a specific selection, matching the
# authored file (this also clears tabSelected on the previously active
# MORT sheet).
$ws.Activate()
$ws.Range("F19").Select() | Out-Null
